$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3821.0527
$ws.Range("I74").Value = 3700
$ws.Range("J74").Value = 3876.923
$ws.Range("K74").Value = 3700
$ws.Range("L74").Value = 3876.923
$ws.Range("M74").Value = -2764
$ws.Range("N74").Value = -5748.923
$ws.Range("H77").Value = 3821.0527
$ws.Range("I77").Value = 3700
$ws.Range("J77").Value = 3876.923
$ws.Range("K77").Value = 18500
$ws.Range("L77").Value = 19384.615
$ws.Range("M77").Value = -13820
$ws.Range("N77").Value = -28744.615
$ws.Range("H111").Value = 2950
$ws.Range("I111").Value = 2900
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 8700
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -5633
$ws.Range("N111").Value = -15134
$ws.Range("H133").Value = 15474.875
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 15474.875
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 15474.875
$ws.Range("N133").Value = -25594.875
$ws.Range("H135").Value = 1323.8823
$ws.Range("I135").Value = 1323.8823
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11914.9407
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9379.940699999999
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1292.762
$ws.Range("I2").Value = 829.8461
$ws.Range("J2").Value = 2045
$ws.Range("K2").Value = 829.8461
$ws.Range("L2").Value = 2045
$ws.Range("M2").Value = -716.8461
$ws.Range("N2").Value = -2271
$ws.Range("H32").Value = 19471.666
$ws.Range("I32").Value = 1755.5
$ws.Range("J32").Value = 134626.75
$ws.Range("K32").Value = 1755.5
$ws.Range("L32").Value = 134626.75
$ws.Range("M32").Value = -1468.5
$ws.Range("N32").Value = -135200.75
$ws.Range("H45").Value = 1559
$ws.Range("I45").Value = 1234.6666
$ws.Range("J45").Value = 2142.8
$ws.Range("K45").Value = 1234.6666
$ws.Range("L45").Value = 2142.8
$ws.Range("M45").Value = -857.6666
$ws.Range("H97").Value = 25649206
$ws.Range("I97").Value = 27786472
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 27786472
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -27785976
$ws.Range("N97").Value = -2992
$ws.Range("H110").Value = 558
$ws.Range("I110").Value = 589.6
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 589.6
$ws.Range("L110").Value = 400
$ws.Range("M110").Value = 1455.4
$ws.Range("N110").Value = -4490
$ws.Range("H116").Value = 1292.762
$ws.Range("I116").Value = 829.8461
$ws.Range("J116").Value = 2045
$ws.Range("K116").Value = 829.8461
$ws.Range("L116").Value = 2045
$ws.Range("M116").Value = 1464.1539
$ws.Range("N116").Value = -6633
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1292.762
$ws.Range("I3").Value = 829.8461
$ws.Range("J3").Value = 2045
$ws.Range("K3").Value = 829.8461
$ws.Range("L3").Value = 2045
$ws.Range("M3").Value = -715.8461
$ws.Range("N3").Value = -2273
$ws.Range("H94").Value = 829.24
$ws.Range("I94").Value = 875.5789
$ws.Range("J94").Value = 682.5
$ws.Range("K94").Value = 875.5789
$ws.Range("L94").Value = 682.5
$ws.Range("M94").Value = -424.5789
$ws.Range("N94").Value = -1584.5
$ws.Range("H99").Value = 2106.7
$ws.Range("I99").Value = 2124.75
$ws.Range("J99").Value = 2094.6667
$ws.Range("K99").Value = 2124.75
$ws.Range("L99").Value = 2094.6667
$ws.Range("M99").Value = -626.75
$ws.Range("N99").Value = -5090.6667
$ws.Range("H105").Value = 3365.4827
$ws.Range("I105").Value = 3136.842
$ws.Range("J105").Value = 3799.9
$ws.Range("K105").Value = 3136.842
$ws.Range("L105").Value = 3799.9
$ws.Range("M105").Value = -1389.842
$ws.Range("N105").Value = -7293.9
$ws.Range("H107").Value = 1406
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1406
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1406
$ws.Range("N107").Value = -5246
$ws.Range("M107").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1244.2424
$ws.Range("I31").Value = 1003.7931
$ws.Range("J31").Value = 2987.5
$ws.Range("K31").Value = 1003.7931
$ws.Range("L31").Value = 2987.5
$ws.Range("M31").Value = -708.7931
$ws.Range("N31").Value = -3577.5
$ws.Range("H34").Value = 1244.2424
$ws.Range("I34").Value = 1003.7931
$ws.Range("J34").Value = 2987.5
$ws.Range("K34").Value = 1003.7931
$ws.Range("L34").Value = 2987.5
$ws.Range("M34").Value = -801.7931
$ws.Range("N34").Value = -3391.5
$ws.Range("H58").Value = 2177.389
$ws.Range("I58").Value = 1525.4
$ws.Range("J58").Value = 2992.375
$ws.Range("K58").Value = 1525.4
$ws.Range("L58").Value = 2992.375
$ws.Range("M58").Value = -1322.4
$ws.Range("N58").Value = -3398.375
$ws.Range("H86").Value = 35716136
$ws.Range("I86").Value = 62501710
$ws.Range("J86").Value = 2031.6666
$ws.Range("K86").Value = 62501710
$ws.Range("L86").Value = 2031.6666
$ws.Range("M86").Value = -62500587
$ws.Range("N86").Value = -4277.6666
$ws.Range("H89").Value = 35716136
$ws.Range("I89").Value = 62501710
$ws.Range("J89").Value = 2031.6666
$ws.Range("K89").Value = 312508550
$ws.Range("L89").Value = 10158.333
$ws.Range("M89").Value = -312502934
$ws.Range("N89").Value = -21390.333
$ws.Range("H99").Value = 7813797
$ws.Range("I99").Value = 15625981
$ws.Range("J99").Value = 1612.5
$ws.Range("K99").Value = 15625981
$ws.Range("L99").Value = 1612.5
$ws.Range("M99").Value = -15624483
$ws.Range("N99").Value = -4608.5
$ws.Range("H126").Value = 7813797
$ws.Range("I126").Value = 15625981
$ws.Range("J126").Value = 1612.5
$ws.Range("K126").Value = 46877943
$ws.Range("L126").Value = 4837.5
$ws.Range("M126").Value = -46875473
$ws.Range("N126").Value = -9777.5
$ws.Range("H136").Value = 2177.389
$ws.Range("I136").Value = 1525.4
$ws.Range("J136").Value = 2992.375
$ws.Range("K136").Value = 4576.200000000001
$ws.Range("L136").Value = 8977.125
$ws.Range("M136").Value = -2026.200000000001
$ws.Range("N136").Value = -14077.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 475.53845
$ws.Range("I33").Value = 441.57144
$ws.Range("J33").Value = 515.1667
$ws.Range("K33").Value = 2649.42864
$ws.Range("L33").Value = 3091.0002
$ws.Range("M33").Value = -2366.42864
$ws.Range("N33").Value = -3657.0002
$ws.Range("H131").Value = 1329.9493
$ws.Range("I131").Value = 317.26666
$ws.Range("J131").Value = 1567.2969
$ws.Range("K131").Value = 951.79998
$ws.Range("L131").Value = 4701.8907
$ws.Range("M131").Value = 4088.20002
$ws.Range("N131").Value = -14781.8907
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1145.25
$ws.Range("I22").Value = 467
$ws.Range("J22").Value = 1371.3334
$ws.Range("K22").Value = 467
$ws.Range("L22").Value = 1371.3334
$ws.Range("M22").Value = -172
$ws.Range("N22").Value = -1961.3334
$ws.Range("H27").Value = 1145.25
$ws.Range("I27").Value = 467
$ws.Range("J27").Value = 1371.3334
$ws.Range("K27").Value = 467
$ws.Range("L27").Value = 1371.3334
$ws.Range("M27").Value = -360
$ws.Range("N27").Value = -1585.3334
$ws.Range("H46").Value = 2486
$ws.Range("I46").Value = 3475
$ws.Range("J46").Value = 1167.3334
$ws.Range("K46").Value = 3475
$ws.Range("L46").Value = 1167.3334
$ws.Range("M46").Value = -3287
$ws.Range("N46").Value = -1543.3334
$ws.Range("H55").Value = 379
$ws.Range("I55").Value = 330.7
$ws.Range("J55").Value = 432.66666
$ws.Range("K55").Value = 330.7
$ws.Range("L55").Value = 432.66666
$ws.Range("M55").Value = -157.7
$ws.Range("N55").Value = -778.66666
$ws.Range("H132").Value = 3797
$ws.Range("I132").Value = 3022.2222
$ws.Range("J132").Value = 5437.706
$ws.Range("K132").Value = 9066.6666
$ws.Range("L132").Value = 16313.118
$ws.Range("M132").Value = -6536.6666
$ws.Range("N132").Value = -21373.118
$ws.Range("H136").Value = 5366.3438
$ws.Range("I136").Value = 3656.238
$ws.Range("J136").Value = 8631.091
$ws.Range("K136").Value = 10968.714
$ws.Range("L136").Value = 25893.273
$ws.Range("M136").Value = -8418.714
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 575120.6
$ws.Range("I81").Value = 2001162.1
$ws.Range("J81").Value = 4704
$ws.Range("K81").Value = 4002324.2
$ws.Range("L81").Value = 9408
$ws.Range("M81").Value = -4001263.2
$ws.Range("N81").Value = -11530
$ws.Range("H84").Value = 575120.6
$ws.Range("I84").Value = 2001162.1
$ws.Range("J84").Value = 4704
$ws.Range("K84").Value = 20011621
$ws.Range("L84").Value = 47040
$ws.Range("M84").Value = -20006317
$ws.Range("N84").Value = -57648
$ws.Range("H132").Value = 38465430
$ws.Range("I132").Value = 71433380
$ws.Range("J132").Value = 2835.5
$ws.Range("K132").Value = 214300140
$ws.Range("L132").Value = 8506.5
$ws.Range("M132").Value = -214297610
$ws.Range("N132").Value = -13566.5
$ws.Range("H135").Value = 102518.125
$ws.Range("I135").Value = 30000
$ws.Range("J135").Value = 112877.86
$ws.Range("K135").Value = 30000
$ws.Range("L135").Value = 112877.86
$ws.Range("M135").Value = -24930
$ws.Range("N135").Value = -123017.86
$ws.Range("H136").Value = 8359888
$ws.Range("I136").Value = 12383833
$ws.Range("J136").Value = 2462.077
$ws.Range("K136").Value = 37151499
$ws.Range("L136").Value = 7386.231000000001
$ws.Range("M136").Value = -37148949
$ws.Range("N136").Value = -12486.231
